$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repulled data: update the dSF (column F) values to reflect the refreshed source data.
$updates = @{
    2  = -3
    3  = -9
    6  = -5
    8  = 4
    15 = -5
    20 = -2
    24 = -6
    26 = -3
    28 = -7
    31 = -3
    32 = 3
    33 = -1
    34 = -2
    35 = -1
    36 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
